$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Items")

# Target values for ambush_chance/ambush_resistance/counter_chance/counter_resistance
# (columns BN, BO, BP, BQ) for trinket rows 2-35, per the commit's data re-import.
$targets = @{
    2  = 0.05
    3  = 0.08
    4  = 0.1
    5  = 0.13
    6  = 0.16
    7  = 0.19
    8  = 0.21
    9  = 0.24
    10 = 0.27
    11 = 0.3
    12 = 0.32
    13 = 0.35
    14 = 0.38
    15 = 0.4
    16 = 0.43
    17 = 0.46
    18 = 0.49
    19 = 0.51
    20 = 0.54
    21 = 0.57
    22 = 0.6
    23 = 0.62
    24 = 0.65
    25 = 0.68
    26 = 0.7
    27 = 0.73
    28 = 0.76
    29 = 0.79
    30 = 0.81
    31 = 0.84
    32 = 0.87
    33 = 0.9
    34 = 0.92
    35 = 0.95
}

foreach ($r in $targets.Keys) {
    $v = $targets[$r]
    $ws.Range("BN$r").Value = $v  # ambush_chance
    $ws.Range("BO$r").Value = $v  # ambush_resistance
    $ws.Range("BP$r").Value = $v  # counter_chance
    $ws.Range("BQ$r").Value = $v  # counter_resistance
}
